# Fix To SPP Creation adding first the Consents
# Append 13 new rows (164-176) to column A, each containing "null",
# matching the existing pattern used in rows 4-163.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 164; $r -le 176; $r++) {
    $ws.Cells.Item($r, 1).Value = "null"
}
